# Updated the tournaments table in content
#
# - Fills in the (previously empty) "platform" column (E) for the reward
#   rows that apply to every platform with the literal "%all%", using the
#   same bordered / left+vertically-centered formatting as its neighbours.
# - Switches the active/selected sheet back to "tournaments" (it had been
#   left on "quests") and updates the remembered selection there.

$wb = $excel.ActiveWorkbook

$tournaments = $wb.Worksheets.Item("tournaments")

# --- Cell content: stamp "%all%" into column E for the relevant reward rows ---
$rows = @(41, 42, 43, 71, 72, 73, 74, 75, 76, 97, 98, 99)
foreach ($r in $rows) {
    $cell = $tournaments.Cells.Item($r, 5)   # column E
    $cell.Value = "%all%"
    $cell.HorizontalAlignment = -4131        # xlLeft
    $cell.VerticalAlignment = -4108          # xlCenter
}

# --- Active tab / selection: make "tournaments" the active sheet again ---
[void]$tournaments.Select()
[void]$tournaments.Range("E11").Select()
